$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) INCLUDEPICTURE field paragraph: swap the two embedded file paths
#    (docx_templates <-> invoices) and resize the VML picture shape from
#    103.8pt x 34.9pt to 103.5pt x 35.25pt.
#    This paragraph is a nested field ( begin / instrText / separate /
#    begin / instrText / instrText / instrText / separate / pict / end /
#    end ) whose instrText runs aren't reachable via Find, so we replace the
#    whole paragraph's content with the exact target OOXML.
# ---------------------------------------------------------------------------
$picturePara = $d.Paragraphs(2).Range

$picXml = '<w:p w:rsidR="00C10012" w:rsidRDefault="009B6994"><w:pPr><w:framePr w:wrap="none" w:vAnchor="page" w:hAnchor="page" w:x="1706" w:y="1873"/><w:rPr><w:sz w:val="2"/><w:szCs w:val="2"/></w:rPr></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> INCLUDEPICTURE  "I:\\MAINDRIVE\\Work\\multiply-file\\docx_templates\\media\\image1.png" \* MERGEFORMATINET </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00E75BB6"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="00E75BB6"><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="00E75BB6"><w:instrText>INCLUDEPICTURE  "I:\\MAINDRIVE\\Work\\multiply-file\\invoices\\media\\image1.png" \* MERGEFORMATINET</w:instrText></w:r><w:r w:rsidR="00E75BB6"><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="00E75BB6"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00E75BB6"><w:pict><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:103.5pt;height:35.25pt"><v:imagedata r:id="rId6" r:href="rId7"/></v:shape></w:pict></w:r><w:r w:rsidR="00E75BB6"><w:fldChar w:fldCharType="end"/></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p>'

$picPackage = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><w:body>' + $picXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$picturePara.InsertXML($picPackage) | Out-Null

# ---------------------------------------------------------------------------
# 2) ACCOUNT NUMBER value: turn the literal "252265874458" into the
#    "{{account_number}}" merge placeholder, split across two runs so the
#    relocated "_GoBack" bookmark can sit between them.
# ---------------------------------------------------------------------------
$accountRange = $d.Content
$found = $accountRange.Find.Execute("252265874458", $true, $false, $false, $false, $false, $true, 1, $false, "{{account_number}}", 2)

if ($found) {
    $splitPoint = $accountRange.Start + "{{account_number".Length
    $bookmarkRange = $d.Range($splitPoint, $splitPoint)
    # Bookmarks.Add moves the bookmark if "_GoBack" already exists elsewhere,
    # which also takes care of removing it from its old location (right
    # before the "25/04/2021 To 25/05/2021" period paragraph).
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)
}

Write-Output "done"
